$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B23").Value = 6333
$ws.Range("C23").Value = 997
$ws.Range("D23").Value = 5901667
$ws.Range("E23").Value = 931.8912048002527
$ws.Range("F23").Value = 8.665065202470835
$ws.Range("G23").Value = 3.746097814776284
$ws.Range("H23").Value = 26.45761770146562
